$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.804.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.07%  "

$ws.Range("D3").Value = "'3.571.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'575.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.59%  "

$ws.Range("D6").Value = "'187.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.29%  "

$ws.Range("D7").Value = "'0.632"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.94%  "

$ws.Range("D8").Value = "'3.567.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("E10").Value = "  -3.55%  "

$ws.Range("D11").Value = "'0.657"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.50%  "

$ws.Range("D12").Value = "'55.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.33%  "

$ws.Range("D13").Value = "'0.0000299"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'9.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("D15").Value = "'4.129.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.65%  "

$ws.Range("D16").Value = "'19.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "

$ws.Range("D17").Value = "'3.561.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.57%  "

$ws.Range("D18").Value = "'69.570.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("D19").Value = "'12.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("E21").Value = "  -2.28%  "

$ws.Range("D22").Value = "'471.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.46%  "

$ws.Range("D23").Value = "'19.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.97%  "

$ws.Range("D24").Value = "'5.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.67%  "

$ws.Range("D25").Value = "'4.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.65%  "

$ws.Range("D26").Value = "'88.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.57%  "

$ws.Range("E27").Value = "  -2.48%  "

$ws.Range("D28").Value = "'10.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.39%  "

$ws.Range("D29").Value = "'9.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").Value = "'32.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("E32").Value = "  +2.07%  "

$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").Value = "'65.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").Value = "'571.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.49%  "

$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("D37").Value = "'0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("E38").Value = "  -4.41%  "

$ws.Range("D39").Value = "'0.396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("E40").Value = "  -6.43%  "

$ws.Range("D41").Value = "'3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.23%  "

$ws.Range("E42").Value = "  +5.15%  "

$ws.Range("D43").Value = "'3.213.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.59%  "

$ws.Range("D44").Value = "'3.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.02%  "

$ws.Range("D45").Value = "'3.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.29%  "

$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("D47").Value = "'9.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.58%  "

$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("D50").Value = "'0.997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("E51").Value = "  -3.47%  "
